# Fill in the 9x9 Sudoku-style grid (B2:J10) with its solved values.
# Row 1 (column headers C1..C9) and column A (row headers R1..R9) already
# contain their labels and are left untouched. Only the previously blank/
# partially-filled interior cells receive values here. The bottom-right
# cell (J10) stays blank, matching the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J")

$grid = @(
    @(1, 2, 3, 4, 5, 6, 7, 8, 9),
    @(4, 5, 6, 7, 8, 9, 1, 2, 3),
    @(7, 8, 9, 1, 2, 3, 4, 5, 6),
    @(2, 3, 4, 5, 6, 7, 8, 9, 1),
    @(5, 6, 7, 8, 9, 1, 2, 3, 4),
    @(8, 9, 1, 2, 3, 4, 5, 6, 7),
    @(3, 4, 5, 6, 7, 8, 9, 0, 1),
    @(6, 7, 8, 9, 0, 1, 3, 4, 5),
    @(9, 0, 1, 3, 4, 5, 6, 7)
)

for ($r = 0; $r -lt $grid.Length; $r++) {
    $rowNum = $r + 2
    $rowVals = $grid[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        # Skip the final, intentionally-blank bottom-right cell (J10).
        if ($rowNum -eq 10 -and $cols[$c] -eq "J") {
            continue
        }
        $cell = $ws.Range($cols[$c] + $rowNum)
        $cell.Value = $rowVals[$c]
        $cell.NumberFormat = "#,##0"
    }
}
